$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 9278060
$ws.Range("I43").Value = 50400.5
$ws.Range("J43").Value = 13891889
$ws.Range("K43").Value = 50400.5
$ws.Range("L43").Value = 13891889
$ws.Range("M43").Value = -50331.5
$ws.Range("N43").Value = -13892027
$ws.Range("H80").Value = 572.7222
$ws.Range("I80").Value = 560.6667
$ws.Range("J80").Value = 575.13336
$ws.Range("K80").Value = 1682.0001
$ws.Range("L80").Value = 1725.40008
$ws.Range("M80").Value = -684.0001
$ws.Range("N80").Value = -3721.40008
$ws.Range("H83").Value = 572.7222
$ws.Range("I83").Value = 560.6667
$ws.Range("J83").Value = 575.13336
$ws.Range("K83").Value = 5046.0003
$ws.Range("L83").Value = 5176.20024
$ws.Range("M83").Value = -54.0002999999997
$ws.Range("N83").Value = -15160.20024
$ws.Range("H132").Value = 5853020
$ws.Range("I132").Value = 6670700.5
$ws.Range("K132").Value = 20012101.5
$ws.Range("M132").Value = -20009571.5
$ws.Range("H138").Value = 589229.6
$ws.Range("I138").Value = 1615
$ws.Range("J138").Value = 714459
$ws.Range("K138").Value = 4845
$ws.Range("L138").Value = 2143377
$ws.Range("M138").Value = 295
$ws.Range("N138").Value = -2153657

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 3009
$ws.Range("J17").Value = 3009
$ws.Range("L17").Value = 3009
$ws.Range("N17").Value = -3355
$ws.Range("H74").Value = 3416
$ws.Range("I74").Value = 2690
$ws.Range("K74").Value = 2690
$ws.Range("M74").Value = -1816
$ws.Range("H77").Value = 3416
$ws.Range("I77").Value = 2690
$ws.Range("K77").Value = 13450
$ws.Range("M77").Value = -9082
$ws.Range("H97").Value = 674.7692
$ws.Range("I97").Value = 527.7778
$ws.Range("J97").Value = 1005.5
$ws.Range("K97").Value = 527.7778
$ws.Range("L97").Value = 1005.5
$ws.Range("M97").Value = -31.77779999999996
$ws.Range("N97").Value = -1997.5
$ws.Range("H102").Value = 16669239
$ws.Range("I102").Value = 20835836
$ws.Range("J102").Value = 2855
$ws.Range("K102").Value = 20835836
$ws.Range("L102").Value = 2855
$ws.Range("M102").Value = -20834214
$ws.Range("N102").Value = -6099
$ws.Range("H106").Value = 21185
$ws.Range("J106").Value = 21185
$ws.Range("L106").Value = 21185
$ws.Range("N106").Value = -23709
$ws.Range("H122").Value = 2885.2
$ws.Range("I122").Value = 2637.3333
$ws.Range("K122").Value = 7911.999899999999
$ws.Range("M122").Value = -5461.999899999999
$ws.Range("H132").Value = 2791.484
$ws.Range("I132").Value = 2414.1365
$ws.Range("K132").Value = 7242.4095
$ws.Range("M132").Value = -4712.4095

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 27778176
$ws.Range("I94").Value = 35714612
$ws.Range("K94").Value = 35714612
$ws.Range("M94").Value = -35714161
$ws.Range("H105").Value = 250001500
$ws.Range("I105").Value = 250001500
$ws.Range("K105").Value = 250001500
$ws.Range("M105").Value = -249999753

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1417.8485
$ws.Range("I31").Value = 1471.8182
$ws.Range("J31").Value = 1390.8636
$ws.Range("K31").Value = 1471.8182
$ws.Range("L31").Value = 1390.8636
$ws.Range("M31").Value = -1176.8182
$ws.Range("N31").Value = -1980.8636
$ws.Range("H34").Value = 1417.8485
$ws.Range("I34").Value = 1471.8182
$ws.Range("J34").Value = 1390.8636
$ws.Range("K34").Value = 1471.8182
$ws.Range("L34").Value = 1390.8636
$ws.Range("M34").Value = -1269.8182
$ws.Range("N34").Value = -1794.8636
$ws.Range("H58").Value = 1734.15
$ws.Range("I58").Value = 1514.8667
$ws.Range("J58").Value = 2392
$ws.Range("K58").Value = 1514.8667
$ws.Range("L58").Value = 2392
$ws.Range("M58").Value = -1311.8667
$ws.Range("N58").Value = -2798
$ws.Range("H105").Value = 846.36365
$ws.Range("I105").Value = 756.6667
$ws.Range("K105").Value = 756.6667
$ws.Range("M105").Value = 990.3333
$ws.Range("H122").Value = 783
$ws.Range("I122").Value = 807.0833
$ws.Range("K122").Value = 2421.2499
$ws.Range("M122").Value = 28.7501000000002
$ws.Range("H132").Value = 2307.125
$ws.Range("I132").Value = 1628.7273
$ws.Range("K132").Value = 4886.1819
$ws.Range("M132").Value = -2356.1819
$ws.Range("H136").Value = 1734.15
$ws.Range("I136").Value = 1514.8667
$ws.Range("J136").Value = 2392
$ws.Range("K136").Value = 4544.6001
$ws.Range("L136").Value = 7176
$ws.Range("M136").Value = -1994.6001
$ws.Range("N136").Value = -12276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 551
$ws.Range("I32").Value = 551
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1653
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1370
$ws.Range("N32").ClearContents()
$ws.Range("H98").Value = 1249
$ws.Range("I98").Value = 1370.5555
$ws.Range("J98").Value = 702
$ws.Range("K98").Value = 4111.666499999999
$ws.Range("L98").Value = 2106
$ws.Range("M98").Value = -2613.666499999999
$ws.Range("N98").Value = -5102
$ws.Range("H122").Value = 1092.6842
$ws.Range("J122").Value = 1122.4117
$ws.Range("L122").Value = 10101.7053
$ws.Range("N122").Value = -15001.7053
$ws.Range("H132").Value = 824.3889
$ws.Range("J132").Value = 882.5
$ws.Range("L132").Value = 7942.5
$ws.Range("N132").Value = -13002.5
$ws.Range("H136").Value = 2078.8125
$ws.Range("I136").Value = 1917.8182
$ws.Range("J136").Value = 2433
$ws.Range("K136").Value = 5753.4546
$ws.Range("L136").Value = 7299
$ws.Range("M136").Value = -653.4546
$ws.Range("N136").Value = -17499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1262.1333
$ws.Range("I102").Value = 1278.2222
$ws.Range("K102").Value = 1278.2222
$ws.Range("M102").Value = 343.7778000000001
$ws.Range("H126").Value = 2161.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1934.8334
$ws.Range("I7").Value = 1840.8
$ws.Range("K7").Value = 1840.8
$ws.Range("M7").Value = -1728.8
$ws.Range("H93").Value = 1001
$ws.Range("J93").Value = 1004
$ws.Range("L93").Value = 1004
$ws.Range("N93").Value = -3500
$ws.Range("H122").Value = 13890557
$ws.Range("I122").Value = 19232224
$ws.Range("K122").Value = 57696672
$ws.Range("M122").Value = -57694222
$ws.Range("H126").Value = 1934.8334
$ws.Range("I126").Value = 1840.8
$ws.Range("K126").Value = 5522.4
$ws.Range("M126").Value = -3052.4
$ws.Range("H132").Value = 2880.1155
$ws.Range("I132").Value = 3029.5
$ws.Range("K132").Value = 9088.5
$ws.Range("M132").Value = -6558.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 50006900
$ws.Range("I62").Value = 83341660
$ws.Range("J62").Value = 4749
$ws.Range("K62").Value = 83341660
$ws.Range("L62").Value = 4749
$ws.Range("M62").Value = -83341036
$ws.Range("N62").Value = -5997
$ws.Range("H65").Value = 50006900
$ws.Range("I65").Value = 83341660
$ws.Range("J65").Value = 4749
$ws.Range("K65").Value = 416708300
$ws.Range("L65").Value = 23745
$ws.Range("M65").Value = -416705180
$ws.Range("N65").Value = -29985
$ws.Range("H132").Value = 3365.742
$ws.Range("I132").Value = 3621.5217
$ws.Range("J132").Value = 2630.375
$ws.Range("K132").Value = 10864.5651
$ws.Range("L132").Value = 7891.125
$ws.Range("M132").Value = -8334.5651
$ws.Range("N132").Value = -12951.125
